$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.549.96'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '1.472.74'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9617'
$ws.Range("E5").Value = '  +5.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '277.44'
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("E7").Value = '  -0.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3086'
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.087'
$ws.Range("E9").Value = '  +6.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.44'
$ws.Range("E10").Value = '  +1.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06632'
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.482'
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.17'
$ws.Range("E14").Value = '  +3.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.164'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9600'
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001024'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '1.471.68'
$ws.Range("E18").Value = '  +2.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05968'
$ws.Range("E19").Value = '  +6.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.91'
$ws.Range("E20").Value = '  +1.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.483'
$ws.Range("E21").Value = '  +1.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.58'
$ws.Range("E22").Value = '  +2.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.26'
$ws.Range("E23").Value = '  +3.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.270'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("D25").Value = '20.553.68'
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.96'
$ws.Range("E26").Value = '  +3.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.112'
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.15'
$ws.Range("E28").Value = '  +1.57%  '
$ws.Range("D29").Value = '1.632.38'
$ws.Range("E29").Value = '  +2.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.84'
$ws.Range("E30").Value = '  +3.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.883'
$ws.Range("E31").Value = '  +1.34%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08002'
$ws.Range("E32").Value = '  +4.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.936'
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.8064'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.226'
$ws.Range("E35").Value = '  +8.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.471'
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05796'
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.715'
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02049'
$ws.Range("E39").Value = '  +3.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9615'
$ws.Range("E40").Value = '  +3.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.39'
$ws.Range("E41").Value = '  +2.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1876'
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.392'
$ws.Range("E43").Value = '  +2.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5270'
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.517'
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.10'
$ws.Range("E46").Value = '  +2.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.11'
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5199'
$ws.Range("E48").Value = '  +2.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.817'
$ws.Range("E49").Value = '  +4.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06449'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9909'
$ws.Range("E51").Value = '  +0.21%  '
